$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B7").Value = "SingleUseId4"
$ws.Range("C7").Value = "Default"
$ws.Range("D7").Value = "Left"
$ws.Range("E7").Value = "LTR"
$ws.Range("F7").Value = "Settings Menu"
